# Deliverables Tracking.xlsx - apply "Updated Professionalism and HDP deliverables" edit
$wb = $excel.ActiveWorkbook

$wsProf = $wb.Worksheets.Item("Professionalism")
$wsHdp  = $wb.Worksheets.Item("Hardware Development Process")

# ---------------------------------------------------------------------------
# Professionalism sheet
# ---------------------------------------------------------------------------

# The "Portoflio" column (G) is no longer tracked on this sheet
$wsProf.Columns.Item(7).Delete()

# Give the task rows real scheduling data instead of being bare text rows.
# Re-use the existing "left-aligned, wrapped" look but with the lighter
# (non-bold) 13pt Arial font used elsewhere for these task rows.
$wsProf.Range("A6:A8").Font.Name = "Arial"
$wsProf.Range("A6:A8").Font.Size = 13
$wsProf.Range("A6:A8").Font.Bold = $false

# Pull the date number-format already used on the HDP sheet so the new date
# cells look consistent across the workbook.
$wsHdp.Range("D6:E6").Copy()
$wsProf.Range("C6:D7").PasteSpecial(-4122)  # xlPasteFormats

$wsProf.Range("A6").Value = "Review, Release, Version Control Lists"
$wsProf.Range("B6").Value = 1
$wsProf.Range("C6").Value = 42921
$wsProf.Range("D6").Value = 42935
$wsProf.Range("E6").Value = "6:00pm"
$wsProf.Rows.Item(6).RowHeight = 32

$wsProf.Range("A7").Value = "Choose a task list format"
$wsProf.Range("B7").Value = 1
$wsProf.Range("C7").Value = 42921
$wsProf.Range("D7").Value = 42935
$wsProf.Range("E7").Value = "6:00pm"
$wsProf.Rows.Item(7).RowHeight = 16

# The old "Task enumeration for Product Definition" row is now blank
$wsProf.Range("A8").ClearContents()
$wsProf.Rows.Item(8).RowHeight = 16

$wsProf.PageSetup.Orientation = 1  # xlPortrait

# Session date moved from a placeholder to the actual term
$wsProf.Range("B2").Value = "Summer 2017"

# ---------------------------------------------------------------------------
# Hardware Development Process sheet
# ---------------------------------------------------------------------------

# Add the term label next to the existing "Independent" note
$wsHdp.Range("D2").Value = "Summer 2017"

# Eagle Layout now runs a bit longer, and "End of Day" became "Beginning of
# Class" as the due-time label for these two deliverables.
$wsHdp.Range("E20").Value = 42924
$wsHdp.Range("F20").Value = "Beginning of Class"

$wsHdp.Range("E21").Value = 42938
$wsHdp.Range("F21").Value = "Beginning of Class"

# Cable Definition isn't applicable for this project, so its schedule is
# replaced with "NA" placeholders.
$wsHdp.Range("C22:F22").ClearContents()
$wsHdp.Range("C22:F22").Value = "NA"
$wsHdp.Rows.Item(22).RowHeight = 16

# ---------------------------------------------------------------------------
# Active sheet / selections
# ---------------------------------------------------------------------------
$wsProf.Activate()
$wsProf.Range("A6:E7").Select()

$wsHdp.Range("F20").Select()
